# Applies the "Automatic update of files." commit to the "Artfynd" sheet:
# rows 14-21 get their observation data permuted (each destination row's
# A/B/D/E/F/G/H/I/M fields, and for rows 16/19 also Y/AA, are replaced by
# another source row's values). Columns C, K, L, N, O, P..X, Z, AB..AY are
# identical across these rows already, so they are left untouched - and so
# is any destination field that happens to already hold the correct value
# (writing it again would be a harmless no-op content-wise, but re-typing
# a cell that doesn't change its value is avoided here to keep the edit
# minimal/surgical, matching the original diff exactly).
#
# A/B/E are genuine numeric cells (Id / Taxonsorteringsordning / TaxonId)
# and are written as numbers. D/F/G/H/I/M/Y/AA are text cells in the
# source file (even "1" and the ISO-looking dates are stored as text, not
# a number/date), so they are written with a leading apostrophe to force
# text entry the way Excel's COM layer would for a user typing into a
# text cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $text) {
    $ws.Range($addr).Value = "'" + $text
}

# Target state for each destination row (source row noted in comment).
$rows = @{
  14 = @{ A=112281199; B=57103;          E=103057; F="Sävsparv";      G="Emberiza schoeniclus";                                       I="";  M="" }                     # <- was row 19
  15 = @{ A=112279542; B=56841; D="NT";  E=103001; F="Rödvingetrast"; G="Turdus iliacus";                H="Linnaeus, 1766";                  M="" }                     # <- was row 17
  16 = @{ A=112281154; B=56575;          E=103021; F="Talltita";      G="Poecile montanus";              H="(Conrad von Baldenstein, 1827)";  I="";  M="födosökande" }   # <- was row 18
  17 = @{ A=112279543; B=56847;          E=102999; F="Björktrast";    G="Turdus pilaris";                H="Linnaeus, 1758" }                                           # <- was row 20
  18 = @{ A=112281210; B=57076; D="VU";  E=103053; F="Lappsparv";     G="Calcarius lapponicus";          H="(Linnaeus, 1758)";               I="1" }                    # <- was row 21
  19 = @{ A=112292314; B=56446;          E=100049; F="Spillkråka";    G="Dryocopus martius";                                                  I="1"; M="lockläte, övriga läten" }  # <- was row 16
  20 = @{ A=112281233; B=57042; D="EN";  E=103042; F="Grönfink";      G="Chloris chloris";               H="(Linnaeus, 1758)";                       M="födosökande" }   # <- was row 15
  21 = @{ A=112279516; B=56321; D="NT";  E=100001; F="Duvhök";        G="Accipiter gentilis" }                                                                          # <- was row 14
}

# Start-/slutdatum only actually change for rows 16 and 19 (they swap with
# each other); every other row in range keeps its existing date, so it is
# not touched here.
$dates = @{
  16 = "2023-09-14"
  19 = "2023-09-15"
}

$numericFields = @("A", "B", "E")
$textFields = @("D", "F", "G", "H", "I", "M")

foreach ($r in $rows.Keys) {
    $row = $rows[$r]
    foreach ($f in $numericFields) {
        if ($row.ContainsKey($f)) {
            $ws.Range("$f$r").Value = $row[$f]
        }
    }
    foreach ($f in $textFields) {
        if ($row.ContainsKey($f)) {
            Set-TextCell $ws "$f$r" $row[$f]
        }
    }
}

foreach ($r in $dates.Keys) {
    Set-TextCell $ws "Y$r" $dates[$r]
    Set-TextCell $ws "AA$r" $dates[$r]
}
